# Updates cryptos list values (coin name/link swaps, prices, 1h volume %)
# matching the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.420.01"
$ws.Range("E2").Value = "  +0.63%  "

# Row 3
$ws.Range("D3").Value = "1.607.99"
$ws.Range("E3").Value = "  +0.97%  "

# Row 5
$ws.Range("D5").Value = "'211.73"
$ws.Range("E5").Value = "  -0.68%  "

# Row 6
$ws.Range("D6").Value = "'0.498"
$ws.Range("E6").Value = "  -0.77%  "

# Row 8
$ws.Range("E8").Value = "  -0.67%  "

# Row 9
$ws.Range("E9").Value = "  -0.23%  "

# Row 10
$ws.Range("D10").Value = "'19.24"
$ws.Range("E10").Value = "  +1.57%  "

# Row 11
$ws.Range("E11").Value = "  -0.44%  "

# Row 12
$ws.Range("D12").Value = "1.834.26"
$ws.Range("E12").Value = "  +0.95%  "

# Row 13
$ws.Range("D13").Value = "1.596.96"
$ws.Range("E13").Value = "  +0.27%  "

# Row 14
$ws.Range("E14").Value = "  -0.15%  "

# Row 15
$ws.Range("E15").Value = "  -0.47%  "

# Row 16
$ws.Range("D16").Value = "'63.45"
$ws.Range("E16").Value = "  -0.75%  "

# Row 17: Bitcoin/WrappedBTC <-> BitcoinCash row swap
$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D17").Value = "'234.23"
$ws.Range("E17").Value = "  +9.09%  "

# Row 18: BitcoinCash <-> WrappedBTC row swap
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "26.415.08"
$ws.Range("E18").Value = "  +0.60%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0724"
$ws.Range("E19").Value = "  -0.21%  "

# Row 20
$ws.Range("D20").Value = "'7.59"
$ws.Range("E20").Value = "  +2.98%  "

# Row 21
$ws.Range("E21").Value = "  -0.17%  "

# Row 22
$ws.Range("E22").Value = "  -0.54%  "

# Row 23
$ws.Range("E23").Value = "  +4.67%  "

# Row 24
$ws.Range("D24").Value = "'8.99"
$ws.Range("E24").Value = "  -0.68%  "

# Row 25
$ws.Range("D25").Value = "'146.64"
$ws.Range("E25").Value = "  +1.16%  "

# Row 26
$ws.Range("E26").Value = "  -0.06%  "

# Row 27
$ws.Range("E27").Value = "  -0.06%  "

# Row 28
$ws.Range("E28").Value = "  +0.71%  "

# Row 29
$ws.Range("D29").Value = "'15.44"
$ws.Range("E29").Value = "  +2.10%  "

# Row 30
$ws.Range("E30").Value = "  +0.76%  "

# Row 31
$ws.Range("E31").Value = "  -0.74%  "

# Row 32
$ws.Range("D32").Value = "1.489.37"
$ws.Range("E32").Value = "  +5.16%  "

# Row 33
$ws.Range("D33").Value = "'3.23"
$ws.Range("E33").Value = "  +1.09%  "

# Row 34
$ws.Range("E34").Value = "  -1.30%  "

# Row 35
$ws.Range("E35").Value = "  -0.53%  "

# Row 36
$ws.Range("E36").Value = "  +0.89%  "

# Row 37
$ws.Range("D37").Value = "'0.564"
$ws.Range("E37").Value = "  -2.43%  "

# Row 38
$ws.Range("E38").Value = "  -0.27%  "

# Row 39
$ws.Range("D39").Value = "'0.823"
$ws.Range("E39").Value = "  +0.07%  "

# Row 40
$ws.Range("D40").Value = "'5.79"
$ws.Range("E40").Value = "  +0.16%  "

# Row 41
$ws.Range("E41").Value = "  -0.05%  "

# Row 42
$ws.Range("D42").Value = "'2.19"
$ws.Range("E42").Value = "  +1.28%  "

# Row 43
$ws.Range("D43").Value = "'0.927"
$ws.Range("E43").Value = "  -5.57%  "

# Row 44
$ws.Range("D44").Value = "1.746.71"
$ws.Range("E44").Value = "  +0.98%  "

# Row 45
$ws.Range("D45").Value = "'0.761"
$ws.Range("E45").Value = "  -0.17%  "

# Row 46
$ws.Range("D46").Value = "'60.98"
$ws.Range("E46").Value = "  +0.04%  "

# Row 47
$ws.Range("D47").Value = "'89.76"
$ws.Range("E47").Value = "  +3.10%  "

# Row 48: BabyDogeCoin -> RenderToken row swap
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.48"
$ws.Range("E48").Value = "  -0.70%  "

# Row 49: RenderToken -> Cronos row swap
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.0501"
$ws.Range("E49").Value = "  -1.08%  "

# Row 50: Cronos -> Algorand row swap
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.0958"
$ws.Range("E50").Value = "  +0.55%  "

# Row 51: Algorand -> EnergySwap row swap
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.46"
$ws.Range("E51").Value = "  +0.65%  "
